# Apply the "4th semester" update:
#  - Lecturers: bump Manuell Koschuch's end-date, tweak selection
#  - Datenbank: select the whole used range
#  - IoT: move the selection
#  - Mobile App Development: move the selection, stop being the active tab
#  - New sheet "Secure Admin Tools" (copied from Mobile App Development),
#    with its own two data rows, and it becomes the active tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Lecturers: Manuell Koschuch's "not available" end date moves out,
#    and the remembered selection changes.
# ---------------------------------------------------------------------
$lecturers = $wb.Worksheets.Item("Lecturers")
$lecturers.Cells.Item(5, 3).Value = 45117
$lecturers.Range("E12").Select()

# ---------------------------------------------------------------------
# 2. Datenbank: whole table now selected.
# ---------------------------------------------------------------------
$datenbank = $wb.Worksheets.Item("Datenbank")
$datenbank.Range("A1:I4").Select()

# ---------------------------------------------------------------------
# 3. IoT: selection moved.
# ---------------------------------------------------------------------
$iot = $wb.Worksheets.Item("IoT")
$iot.Range("G8").Select()

# ---------------------------------------------------------------------
# 4. Mobile App Development: selection moved (it stops being the
#    active tab once the new sheet is created/activated below).
# ---------------------------------------------------------------------
$mobileApp = $wb.Worksheets.Item("Mobile App Development")
$mobileApp.Range("G4").Select()

# ---------------------------------------------------------------------
# 5. New sheet "Secure Admin Tools", modelled on "Mobile App Development".
# ---------------------------------------------------------------------
$mobileApp.Copy([System.Reflection.Missing]::Value, $mobileApp)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Secure Admin Tools"

# Row 2
$newSheet.Cells.Item(2, 1).Value = "Secure Admin Tools"
$newSheet.Cells.Item(2, 2).Value = "Lecture"
$newSheet.Cells.Item(2, 3).Value = 20
$newSheet.Cells.Item(2, 4).Value = 2
$newSheet.Cells.Item(2, 5).Value = 2
$newSheet.Cells.Item(2, 6).Value = "on-site"
$newSheet.Cells.Item(2, 7).Value = "Manuell Koschuch, Silvia Schmidt"
$newSheet.Cells.Item(2, 8).Value = 44971
$newSheet.Cells.Item(2, 9).Value = 45117

# Row 3
$newSheet.Cells.Item(3, 1).Value = "Secure Admin Tools"
$newSheet.Cells.Item(3, 2).Value = "Exercise"
$newSheet.Cells.Item(3, 3).Value = 10
$newSheet.Cells.Item(3, 4).Value = 1
$newSheet.Cells.Item(3, 5).Value = 2
$newSheet.Cells.Item(3, 6).Value = "online"
$newSheet.Cells.Item(3, 7).Value = "Manuell Koschuch, Silvia Schmidt"
$newSheet.Cells.Item(3, 8).Value = 44974
$newSheet.Cells.Item(3, 9).Value = 45113

# Old row 4 (carried over from the copied sheet) is no longer needed.
$newSheet.Rows.Item(4).Delete()

# Header cell A1 uses the same "explicit black" bold style as the rest
# of row 1 (copy it from a cell that already has that exact style).
$pm = $wb.Worksheets.Item("ProjektManagement")
$pm.Cells.Item(1, 2).Copy()
$newSheet.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null

# G2/G3 (lecturer names) use the plain "explicit black" style instead
# of the style inherited from the copied sheet.
$students = $wb.Worksheets.Item("Students")
$students.Cells.Item(52, 1).Copy()
$newSheet.Cells.Item(2, 7).PasteSpecial(-4122) | Out-Null
$newSheet.Cells.Item(3, 7).PasteSpecial(-4122) | Out-Null

$newSheet.Columns.Item(1).ColumnWidth = 17.36328125

$newSheet.Range("D7").Select()
$newSheet.Activate()

$excel.CutCopyMode = 0
